$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-7 with the new TPM-derived values, then remove the old rows 8-10
# that are no longer part of the dataset (sending cluster "ECs" rows removed,
# remaining clusters shifted up).

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Sfrp2"
$ws.Range("C2").Value = "Fzd5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 24.874192
$ws.Range("H2").Value = 74.622576
$ws.Range("I2").Value = 0.9919386828123152
$ws.Range("J2").Value = 0.9919386828123153
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.281109666666667
$ws.Range("N2").Value = 9.843329
$ws.Range("O2").Value = 0.2779739143628921
$ws.Range("P2").Value = 0.2779739143628921
$ws.Range("Q2").Value = 81.61495182172267
$ws.Range("R2").Value = 734.534566395504
$ws.Range("S2").Value = 0.2757330784693105
$ws.Range("T2").Value = 0.2757330784693106

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Sfrp2"
$ws.Range("C3").Value = "Fzd5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 24.874192
$ws.Range("H3").Value = 74.622576
$ws.Range("I3").Value = 0.9919386828123152
$ws.Range("J3").Value = 0.9919386828123153
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.153936333333334
$ws.Range("N3").Value = 18.461809
$ws.Range("O3").Value = 0.5213583040808726
$ws.Range("P3").Value = 0.5213583040808725
$ws.Range("Q3").Value = 153.0741939111093
$ws.Range("R3").Value = 1377.667745199984
$ws.Range("S3").Value = 0.5171554694232433
$ws.Range("T3").Value = 0.5171554694232432

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Sfrp2"
$ws.Range("C4").Value = "Fzd5"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 24.874192
$ws.Range("H4").Value = 74.622576
$ws.Range("I4").Value = 0.9919386828123152
$ws.Range("J4").Value = 0.9919386828123153
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.368614333333333
$ws.Range("N4").Value = 7.105843
$ws.Range("O4").Value = 0.2006677815562353
$ws.Range("P4").Value = 0.2006677815562353
$ws.Range("Q4").Value = 58.91736770128533
$ws.Range("R4").Value = 530.256309311568
$ws.Range("S4").Value = 0.1990501349197615
$ws.Range("T4").Value = 0.1990501349197615

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Sfrp2"
$ws.Range("C5").Value = "Fzd5"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2021483333333333
$ws.Range("H5").Value = 0.606445
$ws.Range("I5").Value = 0.008061317187684791
$ws.Range("J5").Value = 0.008061317187684791
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.281109666666667
$ws.Range("N5").Value = 9.843329
$ws.Range("O5").Value = 0.2779739143628921
$ws.Range("P5").Value = 0.2779739143628921
$ws.Range("Q5").Value = 0.6632708506005557
$ws.Range("R5").Value = 5.969437655405001
$ws.Range("S5").Value = 0.002240835893581602
$ws.Range("T5").Value = 0.002240835893581602

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Sfrp2"
$ws.Range("C6").Value = "Fzd5"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2021483333333333
$ws.Range("H6").Value = 0.606445
$ws.Range("I6").Value = 0.008061317187684791
$ws.Range("J6").Value = 0.008061317187684791
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.153936333333334
$ws.Range("N6").Value = 18.461809
$ws.Range("O6").Value = 0.5213583040808726
$ws.Range("P6").Value = 0.5213583040808725
$ws.Range("Q6").Value = 1.244007973222778
$ws.Range("R6").Value = 11.196071759005
$ws.Range("S6").Value = 0.004202834657629332
$ws.Range("T6").Value = 0.004202834657629331

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Sfrp2"
$ws.Range("C7").Value = "Fzd5"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2021483333333333
$ws.Range("H7").Value = 0.606445
$ws.Range("I7").Value = 0.008061317187684791
$ws.Range("J7").Value = 0.008061317187684791
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.368614333333333
$ws.Range("N7").Value = 7.105843
$ws.Range("O7").Value = 0.2006677815562353
$ws.Range("P7").Value = 0.2006677815562353
$ws.Range("Q7").Value = 0.4788114397927778
$ws.Range("R7").Value = 4.309302958135
$ws.Range("S7").Value = 0.001617646636473857
$ws.Range("T7").Value = 0.001617646636473857

# Remove the now-obsolete trailing rows (old rows 8-10 covered the "MuSCs" sending
# cluster before the data shift; that data now lives in rows 5-7, so the leftover
# rows 8-10 must be deleted to shrink the used range down to A1:T7).
$ws.Range("A8:T10").Delete()
